$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $style = $rng.Style
    $rng.Value = "'" + $value
    $rng.Style = $style
}

Set-TextValue "E2" "0.63%"
Set-TextValue "D3" "26.36"
Set-TextValue "E3" "5.61%"
Set-TextValue "D4" "5.085"
Set-TextValue "E4" "1.43%"
Set-TextValue "D5" "0.05601"
Set-TextValue "E5" "-0.16%"
Set-TextValue "D6" "6.480"
Set-TextValue "E6" "-0.83%"
Set-TextValue "D7" "0.8133"
Set-TextValue "E7" "0.37%"
Set-TextValue "D8" "0.8471"
Set-TextValue "E8" "0.96%"
Set-TextValue "D9" "0.02844"
Set-TextValue "E9" "0.02%"
Set-TextValue "D10" "0.09385"
Set-TextValue "E10" "-0.28%"
Set-TextValue "D11" "0.001516"
Set-TextValue "E11" "0.09%"
Set-TextValue "D12" "0.0005991"
Set-TextValue "E12" "0.51%"
Set-TextValue "D13" "0.006217"
Set-TextValue "E13" "-0.30%"
Set-TextValue "D14" "3.608"
Set-TextValue "E14" "3.08%"
Set-TextValue "D15" "3.011"
Set-TextValue "E15" "0.78%"
Set-TextValue "E16" "-1.73%"
Set-TextValue "E17" "0.72%"
Set-TextValue "D18" "0.1337"
Set-TextValue "E18" "-0.08%"
Set-TextValue "D19" "0.07019"
Set-TextValue "E19" "0.94%"
Set-TextValue "D20" "0.03189"
Set-TextValue "E20" "-3.08%"
Set-TextValue "D21" "0.1297"
Set-TextValue "E21" "0.40%"
Set-TextValue "D22" "3.742"
Set-TextValue "E22" "-0.51%"
Set-TextValue "D23" "0.04650"
Set-TextValue "E23" "-0.77%"
Set-TextValue "E24" "-1.45%"
Set-TextValue "D25" "0.001248"
Set-TextValue "E25" "0.37%"
Set-TextValue "D26" "0.004589"
Set-TextValue "E26" "1.49%"
Set-TextValue "E27" "-0.98%"
Set-TextValue "E28" "-0.10%"
Set-TextValue "E29" "--%"
Set-TextValue "E30" "--%"
Set-TextValue "E31" "--%"
Set-TextValue "E32" "--%"
Set-TextValue "E33" "--%"
Set-TextValue "E34" "--%"
Set-TextValue "E35" "--%"
Set-TextValue "E36" "--%"
Set-TextValue "E37" "--%"
Set-TextValue "E38" "--%"
Set-TextValue "E39" "--%"
Set-TextValue "D40" "0.03675"
Set-TextValue "E40" "1.20%"
Set-TextValue "D41" "0.006231"
Set-TextValue "E41" "85.21%"
Set-TextValue "E42" "0.43%"
Set-TextValue "E43" "-8.17%"
Set-TextValue "D44" "0.008776"
Set-TextValue "E44" "5.00%"
Set-TextValue "D45" "0.00005296"
Set-TextValue "E45" "0.55%"
Set-TextValue "E46" "0.00%"
Set-TextValue "E47" "-40.00%"
Set-TextValue "D48" "0.002310"
Set-TextValue "E48" "9.61%"
Set-TextValue "E49" "0.00%"
Set-TextValue "E50" "0.00%"
Set-TextValue "E51" "--%"
